$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.014
$ws.Range("C2").Value = -0.12
$ws.Range("D2").Value = 0.078
$ws.Range("E2").Value = -0.055
$ws.Range("F2").Value = 0.056
$ws.Range("G2").Value = -0.012
$ws.Range("H2").Value = 0.07
$ws.Range("I2").Value = -0.078
$ws.Range("K2").Value = -0.078
$ws.Range("L2").Value = -0.021
$ws.Range("M2").Value = -0.033

$ws.Range("B3").Value = 0.05
$ws.Range("C3").Value = 0.022
$ws.Range("D3").Value = 0.221
$ws.Range("E3").Value = 0.022
$ws.Range("F3").Value = 0.131
$ws.Range("G3").Value = 0.042
$ws.Range("H3").Value = 0.086
$ws.Range("I3").Value = 0.097
$ws.Range("K3").Value = 0.063
$ws.Range("L3").Value = 0.117
$ws.Range("M3").Value = -0.024

$ws.Range("B4").Value = 0.047
$ws.Range("C4").Value = -0.016
$ws.Range("D4").Value = 0.244
$ws.Range("E4").Value = -0.025
$ws.Range("F4").Value = 0.092
$ws.Range("G4").Value = 0.036
$ws.Range("H4").Value = 0.095
$ws.Range("I4").Value = 0.021
$ws.Range("K4").Value = -0.007
$ws.Range("L4").Value = 0.062
$ws.Range("M4").Value = -0.038

$ws.Range("B5").Value = 0.088
$ws.Range("C5").Value = 0.298
$ws.Range("D5").Value = 0.554
$ws.Range("E5").Value = 0.049
$ws.Range("F5").Value = 0.529
$ws.Range("G5").Value = 0.302
$ws.Range("H5").Value = 0.158
$ws.Range("I5").Value = 0.303
$ws.Range("K5").Value = 0.576
$ws.Range("L5").Value = 0.646
$ws.Range("M5").Value = 0.027

$ws.Range("B6").Value = 0.04
$ws.Range("C6").Value = 1.681
$ws.Range("D6").Value = -0.156
$ws.Range("E6").Value = 1.592
$ws.Range("F6").Value = 1.781
$ws.Range("G6").Value = 0.582
$ws.Range("H6").Value = 0.176
$ws.Range("I6").Value = 0.94
$ws.Range("K6").Value = 2.625
$ws.Range("L6").Value = 1.467
$ws.Range("M6").Value = 0.077

$ws.Range("B7").Value = 0.007
$ws.Range("C7").Value = 0.199
$ws.Range("D7").Value = -0.031
$ws.Range("E7").Value = 0.191
$ws.Range("F7").Value = 0.207
$ws.Range("G7").Value = 0.088
$ws.Range("H7").Value = 0.03
$ws.Range("I7").Value = 0.13
$ws.Range("K7").Value = 0.267
$ws.Range("L7").Value = 0.181
$ws.Range("M7").Value = 0.014

$ws.Range("B8").Value = 0.769
$ws.Range("C8").Value = 1.009
$ws.Range("D8").Value = 1.14
$ws.Range("E8").Value = 0.925
$ws.Range("F8").Value = 1.073
$ws.Range("G8").Value = 0.988
$ws.Range("H8").Value = 1.028
$ws.Range("I8").Value = 0.809
$ws.Range("K8").Value = 1.013
$ws.Range("L8").Value = 0.965
$ws.Range("M8").Value = 0.555

$ws.Range("B9").Value = -0.058
$ws.Range("C9").Value = 0.118
$ws.Range("D9").Value = -0.12
$ws.Range("E9").Value = 0.116
$ws.Range("F9").Value = 0.122
$ws.Range("G9").Value = 0.009
$ws.Range("H9").Value = -0.052
$ws.Range("I9").Value = 0.062
$ws.Range("K9").Value = 0.186
$ws.Range("L9").Value = 0.103
$ws.Range("M9").Value = -0.038

$ws.Range("B10").Value = -0.042
$ws.Range("C10").Value = 0.844
$ws.Range("D10").Value = -0.148
$ws.Range("E10").Value = 0.733
$ws.Range("F10").Value = 0.835
$ws.Range("G10").Value = 0.306
$ws.Range("H10").Value = 0.064
$ws.Range("I10").Value = 0.501
$ws.Range("K10").Value = 0.919
$ws.Range("L10").Value = 0.782
$ws.Range("M10").Value = -0.017

$ws.Range("B11").Value = -0.01
$ws.Range("C11").Value = 0.182
$ws.Range("D11").Value = -0.04
$ws.Range("E11").Value = 0.19
$ws.Range("F11").Value = 0.179
$ws.Range("G11").Value = 0.074
$ws.Range("H11").Value = 0.015
$ws.Range("I11").Value = 0.141
$ws.Range("K11").Value = 0.249
$ws.Range("L11").Value = 0.171
$ws.Range("M11").Value = -0.003

$ws.Range("B12").Value = -0.401
$ws.Range("C12").Value = -0.302
$ws.Range("D12").Value = -0.747
$ws.Range("E12").Value = -0.379
$ws.Range("F12").Value = -0.406
$ws.Range("G12").Value = -0.361
$ws.Range("H12").Value = -0.462
$ws.Range("I12").Value = -0.357
$ws.Range("J12").Value = -1
$ws.Range("K12").Value = -0.295
$ws.Range("L12").Value = -0.339
$ws.Range("M12").Value = -0.282

$ws.Range("B13").Value = 0.189
$ws.Range("C13").Value = 0.217
$ws.Range("D13").Value = 0.31
$ws.Range("E13").Value = 0.24
$ws.Range("F13").Value = 0.23
$ws.Range("G13").Value = 0.237
$ws.Range("H13").Value = 0.233
$ws.Range("I13").Value = 0.228
$ws.Range("K13").Value = 0.274
$ws.Range("L13").Value = 0.212
$ws.Range("M13").Value = 0.097

$ws.Range("B14").Value = 0.388
$ws.Range("C14").Value = 0.675
$ws.Range("D14").Value = 0.44
$ws.Range("E14").Value = 0.456
$ws.Range("F14").Value = 0.805
$ws.Range("G14").Value = 0.506
$ws.Range("H14").Value = 0.565
$ws.Range("I14").Value = 0.375
$ws.Range("K14").Value = 0.616
$ws.Range("L14").Value = 0.727
$ws.Range("M14").Value = 0.319

$ws.Range("B15").Value = 0.065
$ws.Range("C15").Value = 0.081
$ws.Range("D15").Value = 0.089
$ws.Range("E15").Value = 0.075
$ws.Range("F15").Value = 0.085
$ws.Range("G15").Value = 0.079
$ws.Range("H15").Value = 0.082
$ws.Range("I15").Value = 0.068
$ws.Range("K15").Value = 0.081
$ws.Range("L15").Value = 0.078
$ws.Range("M15").Value = 0.051

